# DOMA-4452: added meter place for import
# Adds a new "Место установки счетчика" (meter installation place) column
# as column S, mirroring the formatting of column R (the last existing column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Header cell S1 ---
$ws.Range("S1").Value = "Место установки счетчика"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial($xlPasteFormats)

# --- Data rows 2-11: alternating Кухня / Сан. узел, formatted like column R ---
$places = @("Кухня", "Сан. узел", "Кухня", "Сан. узел", "Кухня", "Сан. узел", "Кухня", "Сан. узел", "Кухня", "Сан. узел")

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Range("S$row").Value = $places[$i]
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

# --- Match column R's width for the new column S ---
$rWidth = $ws.Range("R1").ColumnWidth()
$ws.Range("S1").ColumnWidth = $rWidth
